# Applies the edit described by the commit "CCDC dataset script and icdc programs scripts"
# to the ICDC_ProgramsPage workbook:
#   - the "stdCntVal" description cell for the Comparative Oncology Program (D2 on the
#     "programs" sheet) had a stray leading newline removed from its text
#   - the whole data range on "programs" (A1:F4) was given a Text ("@") number format
#     (the pre-existing yellow fill on E1 and the pre-existing word-wrap on D2:D4 are kept)
#   - row heights on "programs" were tightened up (row 2 back to the sheet default, rows
#     3 & 4 reduced from 75 to 45) now that the text no longer needs as much vertical room
#   - an explicit default-width column definition for column F was added
#   - the sheet selection ended up as the whole sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("programs")

# --- 1. Fix the Comparative Oncology Program "stdCntVal" text (drop the leading newline) ---
$ws.Range("D2").Value = "The COP is a core resource for CCR investigators who are interested in the use of comparative cancer models."

# --- 2. Apply a Text number format across the whole used range. Doing this *after* the
#        values are written keeps the numeric index/count cells (A/C/E/F columns) as real
#        numbers instead of turning them into text. Cells that already carried other
#        formatting (E1's yellow fill, D2:D4's wrap text) keep that formatting alongside
#        the new number format. ---
$ws.Range("A1:F4").NumberFormat = "@"

# --- 3. Row heights: row 2 no longer needs the extra height it had before, rows 3 & 4
#        shrink from 75 to 45. ---
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45

# --- 4. Give column F (beyond the original E) an explicit standard width, matching the
#        rest of the sheet's formatting being extended to it. ---
$ws.Columns.Item(6).ColumnWidth = 9.14

# --- 5. Leave the sheet with the whole grid selected (matches the saved selection state). ---
$ws.Cells.Select()

Write-Host "Applied ICDC programs page updates"
